# "Add files via upload" - the uploaded workbook reflects a tester going
# through both sheets and filling in the Approved/Rejected review columns.
#
# GA_gentp (sheet 1):
#   - row 95 (TestScenario_6 / create Document) gets flipped from
#     "Approved" to "Rejected", with a reason noted in column J.
#   - row 103 (TestScenario_7) keeps its "Rejected" status but its
#     ReasonToReject note changes from "checking" to "bccggcbg".
# SF_salestp (sheet 2):
#   - every test-case summary row (2-13) is marked "Approved" in column I.
#   - the final row (14) is marked "Rejected" with a "checking" note in
#     column J, and becomes the sheet left active/selected on save.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GA_gentp")
$ws2 = $wb.Worksheets.Item("SF_salestp")

# --- GA_gentp updates ---------------------------------------------------
$ws1.Range("I95").Value = "Rejected"
$ws1.Range("J103").Value = "bccggcbg"
$ws1.Range("J95").Value = "vghghdghd"

# --- SF_salestp updates --------------------------------------------------
for ($r = 2; $r -le 13; $r++) {
    $ws2.Cells.Item($r, 9).Value = "Approved"
}
$ws2.Cells.Item(14, 9).Value = "Rejected"
$ws2.Cells.Item(14, 10).Value = "checking"

# --- Selection / active sheet state --------------------------------------
# Leave GA_gentp's cursor on the newly edited J95 cell...
$ws1.Range("J95").Select() | Out-Null
# ...but the workbook was saved with SF_salestp as the visible/active tab,
# cursor parked on the last edited cell J14.
$ws2.Activate()
$ws2.Range("J14").Select() | Out-Null
